$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5, columns A-G with new values
$ws.Range("A2").Value = 0.03072753924798493
$ws.Range("B2").Value = 0.007846346145009038
$ws.Range("C2").Value = 0.06582808494567871
$ws.Range("D2").Value = 0.09126154182117911
$ws.Range("E2").Value = 0.1066511746974878
$ws.Range("F2").Value = 0.09835801166434541
$ws.Range("G2").Value = 0.02333569645984184

$ws.Range("A3").Value = 0.01489915666857287
$ws.Range("B3").Value = 0.003310804408225765
$ws.Range("C3").Value = -0.00726838962117827
$ws.Range("D3").Value = 0.04882070088206886
$ws.Range("E3").Value = 0.1011448360198526
$ws.Range("F3").Value = 0.06585462082959219
$ws.Range("G3").Value = 0.01015087145969499

$ws.Range("A4").Value = 0.05485779992694893
$ws.Range("B4").Value = 0.01750636503375207
$ws.Range("C4").Value = 0.07172764837741852
$ws.Range("D4").Value = 0.1178446918715434
$ws.Range("E4").Value = 0.1258494983977992
$ws.Range("F4").Value = 0.121715625181588
$ws.Range("G4").Value = 0.02610132981227498

$ws.Range("A5").Value = 0.004053295892611886
$ws.Range("B5").Value = 0.0007101097404372962
$ws.Range("C5").Value = 0.004137439861266331
$ws.Range("D5").Value = 0.0748104645746592
$ws.Range("E5").Value = 0.1010751379315278
$ws.Range("F5").Value = 0.08598177358308073
$ws.Range("G5").Value = 0.01371640142198511

# Add new row 6 with CNN descriptor
$ws.Range("A6").Value = 0.1342448192655817
$ws.Range("B6").Value = 0.0261367725086431
$ws.Range("C6").Value = -0.0216719675809145
$ws.Range("D6").Value = 0.1605071901277799
$ws.Range("E6").Value = 0.2475178865761927
$ws.Range("F6").Value = 0.1947350922724521
$ws.Range("G6").Value = 0.01953489584559281
$ws.Range("H6").Value = "CNN"
$ws.Range("I6").Value = "kmeans"
